# Auto-generated Excel COM-interop script
# Applies scraped market-price updates (currentAveragePrice* / LevePrice* / LeveProfit*)
# columns H-N) across the ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets, matching the upstream diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1529
$ws.Range("I11").Value = 1529
$ws.Range("K11").Value = 1529
$ws.Range("M11").Value = -1389
$ws.Range("H17").Value = 163648.27
$ws.Range("J17").Value = 167714.47
$ws.Range("L17").Value = 503143.41
$ws.Range("N17").Value = -503479.41
$ws.Range("H41").Value = 740
$ws.Range("I41").Value = 317.42856
$ws.Range("K41").Value = 317.42856
$ws.Range("M41").Value = 122.57144
$ws.Range("H54").Value = 26657.715
$ws.Range("I54").Value = 5000
$ws.Range("K54").Value = 5000
$ws.Range("M54").Value = -4514
$ws.Range("H62").Value = 4224.75
$ws.Range("I62").Value = 4166.3335
$ws.Range("K62").Value = 4166.3335
$ws.Range("M62").Value = -3542.3335
$ws.Range("H65").Value = 4224.75
$ws.Range("I65").Value = 4166.3335
$ws.Range("K65").Value = 20831.6675
$ws.Range("M65").Value = -17711.6675
$ws.Range("H98").Value = 3508.5
$ws.Range("I98").Value = 3445.9167
$ws.Range("J98").Value = 3884
$ws.Range("K98").Value = 3445.9167
$ws.Range("L98").Value = 3884
$ws.Range("M98").Value = -1947.9167
$ws.Range("N98").Value = -6880
$ws.Range("H122").Value = 3508.5
$ws.Range("I122").Value = 3445.9167
$ws.Range("J122").Value = 3884
$ws.Range("K122").Value = 10337.7501
$ws.Range("L122").Value = 11652
$ws.Range("M122").Value = -7887.750100000001
$ws.Range("N122").Value = -16552
$ws.Range("H132").Value = 2180.3044
$ws.Range("I132").Value = 2223.5264
$ws.Range("K132").Value = 6670.5792
$ws.Range("M132").Value = -4140.5792

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 95000
$ws.Range("J24").Value = 95000
$ws.Range("L24").Value = 95000
$ws.Range("N24").Value = -95748
$ws.Range("H74").Value = 21603.098
$ws.Range("I74").Value = 1976.4103
$ws.Range("K74").Value = 1976.4103
$ws.Range("M74").Value = -1102.4103
$ws.Range("H77").Value = 21603.098
$ws.Range("I77").Value = 1976.4103
$ws.Range("K77").Value = 9882.0515
$ws.Range("M77").Value = -5514.0515
$ws.Range("H97").Value = 4097.8887
$ws.Range("I97").Value = 4256.7856
$ws.Range("K97").Value = 4256.7856
$ws.Range("M97").Value = -3760.7856
$ws.Range("H100").Value = 95000
$ws.Range("J100").Value = 95000
$ws.Range("L100").Value = 95000
$ws.Range("N100").Value = -97164
$ws.Range("H104").Value = 30990
$ws.Range("J104").Value = 30990
$ws.Range("L104").Value = 30990
$ws.Range("N104").Value = -37978

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 28589156
$ws.Range("I16").Value = 47620264
$ws.Range("J16").Value = 42499.5
$ws.Range("K16").Value = 47620264
$ws.Range("L16").Value = 42499.5
$ws.Range("M16").Value = -47619977
$ws.Range("N16").Value = -43073.5
$ws.Range("H31").Value = 4883.1113
$ws.Range("I31").Value = 4481.5
$ws.Range("J31").Value = 5119.353
$ws.Range("K31").Value = 4481.5
$ws.Range("L31").Value = 5119.353
$ws.Range("M31").Value = -4186.5
$ws.Range("N31").Value = -5709.353
$ws.Range("H34").Value = 4883.1113
$ws.Range("I34").Value = 4481.5
$ws.Range("J34").Value = 5119.353
$ws.Range("K34").Value = 4481.5
$ws.Range("L34").Value = 5119.353
$ws.Range("M34").Value = -4279.5
$ws.Range("N34").Value = -5523.353
$ws.Range("H107").Value = 1901.8
$ws.Range("I107").Value = 1972.1
$ws.Range("J107").Value = 1761.2
$ws.Range("K107").Value = 1972.1
$ws.Range("L107").Value = 1761.2
$ws.Range("M107").Value = -52.09999999999991
$ws.Range("N107").Value = -5601.2
$ws.Range("H113").Value = 28589156
$ws.Range("I113").Value = 47620264
$ws.Range("J113").Value = 42499.5
$ws.Range("K113").Value = 47620264
$ws.Range("L113").Value = 42499.5
$ws.Range("M113").Value = -47618094
$ws.Range("N113").Value = -46839.5
$ws.Range("H141").Value = 273241.84
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 273241.84
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 273241.84
$ws.Range("N141").Value = -283601.84
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 453.25
$ws.Range("I2").Value = 86.40000000000001
$ws.Range("J2").Value = 715.2857
$ws.Range("K2").Value = 518.4000000000001
$ws.Range("L2").Value = 4291.7142
$ws.Range("M2").Value = -405.4000000000001
$ws.Range("N2").Value = -4517.7142
$ws.Range("H7").Value = 1236.4
$ws.Range("I7").Value = 290
$ws.Range("K7").Value = 870
$ws.Range("M7").Value = -758
$ws.Range("H75").Value = 8930670
$ws.Range("J75").Value = 11907285
$ws.Range("L75").Value = 35721855
$ws.Range("N75").Value = -35723851
$ws.Range("H78").Value = 8930670
$ws.Range("J78").Value = 11907285
$ws.Range("L78").Value = 107165565
$ws.Range("N78").Value = -107175549
$ws.Range("H112").Value = 14671.333
$ws.Range("I112").Value = 3000
$ws.Range("J112").Value = 17005.6
$ws.Range("K112").Value = 9000
$ws.Range("L112").Value = 51016.8
$ws.Range("M112").Value = -7892
$ws.Range("N112").Value = -53232.8
$ws.Range("H122").Value = 5556727
$ws.Range("I122").Value = 6667450
$ws.Range("K122").Value = 60007050
$ws.Range("M122").Value = -60004600

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 9000
$ws.Range("J40").Value = 9000
$ws.Range("L40").Value = 9000
$ws.Range("M40").Value = -9302
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5277
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("H106").Value = 17456.666
$ws.Range("J106").Value = 17456.666
$ws.Range("L106").Value = 17456.666
$ws.Range("N106").Value = -19980.666
$ws.Range("H136").Value = 43641.25
$ws.Range("I136").Value = 54337.895
$ws.Range("K136").Value = 163013.685
$ws.Range("M136").Value = -160463.685
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 83583.336
$ws.Range("J46").Value = 83583.336
$ws.Range("L46").Value = 83583.336
$ws.Range("N46").Value = -84045.336
$ws.Range("H134").Value = 83583.336
$ws.Range("J134").Value = 83583.336
$ws.Range("L134").Value = 250750.008
$ws.Range("N134").Value = -255820.008
$ws.Range("H141").Value = 119998
$ws.Range("J141").Value = 119998
$ws.Range("L141").Value = 119998
$ws.Range("N141").Value = -130358

Write-Output "Applied 170 cell updates across 7 sheets"